$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "ammonia",
    "arsenic",
    "barium",
    "cadmium",
    "chloramine",
    "chromium",
    "copper",
    "flouride",
    "bacteria",
    "viruses",
    "lead",
    "nitrates",
    "nitrites",
    "mercury",
    "perchlorate",
    "radium",
    "selenium",
    "silver",
    "uranium",
    "aluminiumdanger",
    "is_safe"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

[void]$ws.Range("B3").Select()
